# Apply the "West Ham_stats.xlsx" edit:
#  1. Increment the day component of every E-column "Age" value (format
#     "YY-DDD") by 1 on every stats sheet (all sheets except "Matches").
#  2. Split the merged "Playing Time" header on the two sheets that have a
#     stray extra "Unnamed: 4_level_0" column (StandardStats / PlayingTime):
#     F1 becomes its own "Unnamed: 4_level_0" label and the "Playing Time"
#     label + merge move one column right, from F1:I1 to G1:I1.
#  3. Rename the stats sheets to their human-readable, space-separated
#     (and in one case ampersand-joined) names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Bump the Age (column E) day count on every stats sheet.
# ---------------------------------------------------------------------
$statSheetNames = @(
    "StandardStats",
    "ShootingStats",
    "PassingStats",
    "PassTypes",
    "GoalShotCreation",
    "DefensiveActions",
    "Possession",
    "PlayingTime",
    "MiscStats"
)

foreach ($sheetName in $statSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $val = $cell.Value2
        if ($val -ne $null -and $val -ne "") {
            $parts = [string]$val -split "-"
            if ($parts.Count -eq 2) {
                $years = $parts[0]
                $days = [int]$parts[1] + 1
                $cell.Value = "{0}-{1:D3}" -f $years, $days
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Fix the "Playing Time" merged header on StandardStats & PlayingTime.
#    Before: F1="Playing Time" (merged F1:I1), G1/H1/I1 blank.
#    After:  F1="Unnamed: 4_level_0", G1="Playing Time" (merged G1:I1),
#            H1/I1 blank - all four cells keep their original style.
# ---------------------------------------------------------------------
$headerSheetNames = @("StandardStats", "PlayingTime")

foreach ($sheetName in $headerSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F1:I1").UnMerge()
    $ws.Range("F1").Value = "Unnamed: 4_level_0"
    $ws.Range("G1").Value = "Playing Time"
    $ws.Range("G1:I1").Merge()

    # Re-apply the original header styling (bold font, centered/top
    # aligned, thin box border) so the cells land back on the same style
    # record instead of the engine minting fresh ones for the merge.
    $hdr = $ws.Range("F1:I1")
    $hdr.Font.Bold = $true
    $hdr.HorizontalAlignment = -4108
    $hdr.VerticalAlignment = -4160
    $hdr.Borders.LineStyle = 1
    $hdr.Borders.Weight = 2
}

# ---------------------------------------------------------------------
# 3. Rename the stats sheets to their spaced-out display names.
# ---------------------------------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"     = "Shooting Stats"
    "PassingStats"      = "Passing Stats"
    "PassTypes"         = "Pass Types"
    "GoalShotCreation"  = "Goal & Shot Creation"
    "DefensiveActions"  = "Defensive Actions"
    "PlayingTime"       = "Playing Time"
    "MiscStats"         = "Miscellaneous Stats"
}

foreach ($oldName in $renames.Keys) {
    $ws = $wb.Worksheets.Item($oldName)
    $ws.Name = $renames[$oldName]
}
